$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last data row (old row 22) - dataset now ends at row 21
$ws.Rows.Item(22).Delete()

# Update sensor reading columns C:H for data rows 2-21 with newly computed values
$ws.Range("C2").Value = -4.713510870933534
$ws.Range("D2").Value = -6.196972131729122
$ws.Range("E2").Value = -6.711660981178287
$ws.Range("F2").Value = 1.058934926986694
$ws.Range("G2").Value = -1.520596265792847
$ws.Range("H2").Value = -2.644282817840576
$ws.Range("C3").Value = -5.068336345255377
$ws.Range("D3").Value = -2.708849310874939
$ws.Range("E3").Value = -9.788534998893738
$ws.Range("F3").Value = 3.271336078643799
$ws.Range("G3").Value = -1.098488450050354
$ws.Range("H3").Value = -2.443155527114868
$ws.Range("C4").Value = 0.0401994585990876
$ws.Range("D4").Value = -2.849318981170666
$ws.Range("E4").Value = -9.811291694641097
$ws.Range("F4").Value = -2.085494041442871
$ws.Range("G4").Value = -0.6218608021736145
$ws.Range("H4").Value = -1.758833765983582
$ws.Range("C5").Value = -1.265282429754738
$ws.Range("D5").Value = -8.358855545520797
$ws.Range("E5").Value = -3.638464748859389
$ws.Range("F5").Value = -6.340017795562744
$ws.Range("G5").Value = 0.8550586104393005
$ws.Range("H5").Value = 4.873482704162598
$ws.Range("C6").Value = -8.772404074668895
$ws.Range("D6").Value = -9.049420595169067
$ws.Range("E6").Value = 1.087344408035282
$ws.Range("F6").Value = -5.074305057525635
$ws.Range("G6").Value = 1.213025689125061
$ws.Range("H6").Value = -1.407586216926575
$ws.Range("C7").Value = -9.486006855964657
$ws.Range("D7").Value = -7.914790928363797
$ws.Range("E7").Value = 4.333469152450566
$ws.Range("F7").Value = -0.499840497970581
$ws.Range("G7").Value = 1.082911372184753
$ws.Range("H7").Value = -0.4966334402561188
$ws.Range("C8").Value = -6.667140722274777
$ws.Range("D8").Value = -5.866719007492064
$ws.Range("E8").Value = 7.366441488265995
$ws.Range("F8").Value = -1.038776397705078
$ws.Range("G8").Value = 4.970457553863525
$ws.Range("H8").Value = -1.859473824501038
$ws.Range("C9").Value = -3.909939825534814
$ws.Range("D9").Value = -4.199757695198058
$ws.Range("E9").Value = 8.474415659904475
$ws.Range("F9").Value = -1.67560338973999
$ws.Range("G9").Value = 1.536020636558533
$ws.Range("H9").Value = -2.598315238952637
$ws.Range("C10").Value = -1.259068042039859
$ws.Range("D10").Value = -3.636529445648192
$ws.Range("E10").Value = 9.419335365295435
$ws.Range("F10").Value = 2.161241292953491
$ws.Range("G10").Value = 1.872912883758545
$ws.Range("H10").Value = -0.7365507483482361
$ws.Range("C11").Value = 1.691281750798222
$ws.Range("D11").Value = -3.219833493232723
$ws.Range("E11").Value = 15.28036725521086
$ws.Range("F11").Value = 0.3984368443489074
$ws.Range("G11").Value = 0.0755945742130279
$ws.Range("H11").Value = -0.2574796974658966
$ws.Range("C12").Value = 0.7758507728576624
$ws.Range("D12").Value = -2.317984580993662
$ws.Range("E12").Value = 10.35012340545653
$ws.Range("F12").Value = -1.350317597389221
$ws.Range("G12").Value = 0.328340083360672
$ws.Range("H12").Value = 0.2541199326515198
$ws.Range("C13").Value = -0.05869728326797519
$ws.Range("D13").Value = -4.645559132099153
$ws.Range("E13").Value = 8.129511415958403
$ws.Range("F13").Value = 2.818379640579224
$ws.Range("G13").Value = -2.691624879837036
$ws.Range("H13").Value = 2.749962568283081
$ws.Range("C14").Value = -7.922254800796521
$ws.Range("D14").Value = -8.533369660377506
$ws.Range("E14").Value = 4.239331245422359
$ws.Range("F14").Value = 1.895667552947998
$ws.Range("G14").Value = -2.188882827758789
$ws.Range("H14").Value = -0.0742201283574104
$ws.Range("C15").Value = -7.12406146526336
$ws.Range("D15").Value = -9.578974485397339
$ws.Range("E15").Value = 5.120484650135045
$ws.Range("F15").Value = 0.1685988008975982
$ws.Range("G15").Value = -1.522581577301025
$ws.Range("H15").Value = 0.997542917728424
$ws.Range("C16").Value = -4.634187221527098
$ws.Range("D16").Value = -9.333477973937988
$ws.Range("E16").Value = 4.373982667922966
$ws.Range("F16").Value = 0.164170041680336
$ws.Range("G16").Value = -3.341738224029541
$ws.Range("H16").Value = 0.8228355050086975
$ws.Range("C17").Value = -2.660379245877256
$ws.Range("D17").Value = -9.284614562988283
$ws.Range("E17").Value = -3.805972993373912
$ws.Range("F17").Value = 0.7200574278831482
$ws.Range("G17").Value = -2.653293132781982
$ws.Range("H17").Value = -0.8017606139183044
$ws.Range("C18").Value = 0.6678269803524017
$ws.Range("D18").Value = -9.062278509140016
$ws.Range("E18").Value = -10.71978342533112
$ws.Range("F18").Value = -0.3216205537319183
$ws.Range("G18").Value = -0.1701259762048721
$ws.Range("H18").Value = 0.732122004032135
$ws.Range("C19").Value = 2.381343364715557
$ws.Range("D19").Value = -7.552583456039422
$ws.Range("E19").Value = 0.6790638566017475
$ws.Range("F19").Value = -4.379751205444336
$ws.Range("G19").Value = 0.5068654417991638
$ws.Range("H19").Value = 1.012203693389893
$ws.Range("C20").Value = -2.60627746582033
$ws.Range("D20").Value = -5.606535911560057
$ws.Range("E20").Value = 9.005005836486802
$ws.Range("F20").Value = -2.21698260307312
$ws.Range("G20").Value = -1.410029649734497
$ws.Range("H20").Value = -0.8793405294418335
$ws.Range("C21").Value = -6.737108409404774
$ws.Range("D21").Value = -5.243253648281096
$ws.Range("E21").Value = 5.65115070343016
$ws.Range("F21").Value = -0.2820670306682586
$ws.Range("G21").Value = 0.3074179291725158
$ws.Range("H21").Value = -1.758681058883667
